# vipul commit after  VPPersonalPC branch merge
#
# Update the "Portal Range" values on the Credentials sheet and move the
# active selection to reflect where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# Portal Range: 35 -> 101, 334 -> 111
$ws.Range("B10").Value = "101"
$ws.Range("B11").Value = "111"

# Leave the selection on B11 (previously B17)
$ws.Range("B11").Select()
